# Apply the "To Dos" workbook update:
#  - Two completed/obsolete notes on the "30.09." sheet are cleared out
#    (their text no longer needed anywhere, so the strings drop out of the table).
#  - On "Tabelle1" several to-do items are marked done (highlighted with the
#    existing "done" fill style) and new follow-up notes are appended, matching
#    the commit: "Belegen in update-K. funktioniert jetzt ... .find -> .findOne".
#  - The active sheet/selection moves from "Notizen Präsi" back to "Tabelle1".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("30.09.")

# --- "30.09." sheet: remove the two now-resolved notes ---------------------
$ws2.Range("A2:A3").ClearContents()

# --- "Tabelle1": add the new to-do / follow-up notes ------------------------
# (Order matters so the strings land at the same shared-string slots as the
# authored workbook: F6, then A33, A34, A32.)
$ws1.Range("F6").Value = "formvontrolname vorne?"
$ws1.Range("A33").Value = "update Methode hinterlegen"
$ws1.Range("A34").Value = "create new erstellen!"
$ws1.Range("A32").Value = "FormGroup befüllen fixen! Weil es ein Array ist?"

# --- Mark items as done: reuse the highlighted "done" style from A10 -------
$doneStyle = $ws1.Range("A10")
$doneStyle.Copy()
$ws1.Range("E7").PasteSpecial(-4122)
$ws1.Range("A9").PasteSpecial(-4122)
$ws1.Range("E9").PasteSpecial(-4122)
$ws1.Range("A16").PasteSpecial(-4122)
$ws1.Range("A17").PasteSpecial(-4122)
$ws1.Range("A24").PasteSpecial(-4122)
$ws1.Range("A32").PasteSpecial(-4122)

# --- Move the active sheet/selection back to "Tabelle1" --------------------
$ws1.Activate()
$ws1.Range("A17").Select()

Write-Output "edit complete"
